# Slide 2, "Content Placeholder 2" shape, 3rd paragraph:
#   "New banking customers have to be added to the bank database (JSON file) to use this application. "
# becomes three runs:
#   "New banking customers have to be added to the bank database (JSON file "
#   "& MongoDB) "
#   "to use this application. "

$p   = $ppt.ActivePresentation
$s   = $p.Slides.Item(2)
$shp = $s.Shapes.Item(2)
$tr  = $shp.TextFrame.TextRange

# The sentence about new customers is the 3rd paragraph in this placeholder.
$para = $tr.Paragraphs(3, 1)

$originalText = $para.Text
$closeParenPos = $originalText.IndexOf(")") + 1   # 1-based index of the ')' after "JSON file"

$firstPart  = "New banking customers have to be added to the bank database (JSON file & MongoDB) "
$secondPart = "& MongoDB) "

# Replace everything through the ") " with the merged first+second text; this keeps
# the trailing "to use this application. " run (and its formatting) untouched.
$para.Characters(1, $closeParenPos + 1).Text = $firstPart

# Now split the just-inserted tail off into its own run so "& MongoDB) " carries its
# own run properties, matching the authored edit.
$para.Characters($firstPart.Length - $secondPart.Length + 1, $secondPart.Length).Text = $secondPart
